# Seed "n_likes" column on the posts sheet and document it with a comment,
# then leave the workbook focused on the users sheet (matching the
# author's final click-through state).

$wb = $excel.ActiveWorkbook

$wsUsers = $wb.Worksheets.Item(1)
$wsPosts = $wb.Worksheets.Item(3)

# --- posts sheet: add the n_likes column -------------------------------
# Copy the header formatting from the neighboring "privacy" header (E1)
# onto the new F1 header cell, then set its label.
$wsPosts.Range("E1").Copy() | Out-Null
$wsPosts.Range("F1").PasteSpecial(-4122) | Out-Null
$wsPosts.Range("F1").Value = "n_likes"

# Seed a deterministic number of likes (0-20) per post row (2-31), pulled
# from the dummy user accounts defined on the users/creators pages.
$likes = @(9, 9, 11, 9, 5, 9, 8, 7, 7, 4, 1, 3, 2, 4, 6, 3, 8, 1, 4, 0, 0, 0, 0, 0, 0, 0, 2, 1, 8, 10)
for ($i = 0; $i -lt $likes.Length; $i++) {
    $row = 2 + $i
    $wsPosts.Cells.Item($row, 6).Value = $likes[$i]
}

# Document the new column with a comment explaining where the like counts
# come from.
$wsPosts.Range("F1").AddComment("Aaron Becker:" + [char]10 + "There are up to 20 potential likes from users defined on previous page") | Out-Null

# --- window state: leave selection on the posts sheet, then move back to
#     the users sheet (matches the final activeTab / tabSelected state) --
$wsPosts.Activate()
$wsPosts.Range("F2").Select() | Out-Null

$wsUsers.Activate()
$wsUsers.Range("D10").Select() | Out-Null
